# Update Work Week and Social Spending
# (Gabon GDP per Capita series: revise existing years 1950-2008 with new
# estimates, and append newly available years 2009-2016.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New "GDP per Capita" values for every year 1950-2016, keyed by year.
# Values are written with a leading apostrophe so they stay text cells
# (matching the source file, where the Data column is stored as text),
# exactly as the underlying numeric-looking strings were before the edit.
$values = @{
    1950 = "2958"
    1951 = "3049"
    1952 = "3143"
    1953 = "3237"
    1954 = "3335"
    1955 = "3437"
    1956 = "3539"
    1957 = "3642"
    1958 = "3749"
    1959 = "3857"
    1960 = "3982"
    1961 = "4415"
    1962 = "4497"
    1963 = "4599"
    1964 = "4616"
    1965 = "4626"
    1966 = "4761"
    1967 = "4882"
    1968 = "4927"
    1969 = "5252"
    1970 = "5587"
    1971 = "6027"
    1972 = "6559"
    1973 = "6935"
    1974 = "9083"
    1975 = "9985"
    1976 = "12747"
    1977 = "10849"
    1978 = "8013"
    1979 = "8091"
    1980 = "8402"
    1981 = "8625"
    1982 = "8099"
    1983 = "8281"
    1984 = "8612"
    1985 = "8131"
    1986 = "7818"
    1987 = "6320"
    1988 = "6982"
    1989 = "7422"
    1990 = "7646"
    1991 = "8296.57394557838"
    1992 = "8194.49945140058"
    1993 = "8551.05872261632"
    1994 = "9192.5629850228"
    1995 = "9833.98903886481"
    1996 = "10389.2097947514"
    1997 = "11202.8387820701"
    1998 = "11758.1601210366"
    1999 = "10648.0894892769"
    2000 = "10637.7997208308"
    2001 = "11104.5530426954"
    2002 = "11141.4626736169"
    2003 = "11665.8322643403"
    2004 = "12022.5071831266"
    2005 = "12437.9596840368"
    2006 = "12342.7793300667"
    2007 = "13455.6124838169"
    2008 = "13330.2547689156"
    2009 = "13617.8961017338"
    2010 = "14926.6996403551"
    2011 = "16403"
    2012 = "16932"
    2013 = "17532"
    2014 = "18053"
    2015 = "18389"
    2016 = "18413"
}

# Rows 2-60 already hold years 1950-2008; just refresh column E (Data).
for ($year = 1950; $year -le 2008; $year++) {
    $row = $year - 1950 + 2
    $ws.Cells.Item($row, 5).Value = "'" + $values[$year]
}

# Append the newly reported years 2009-2016 as rows 61-68.
for ($year = 2009; $year -le 2016; $year++) {
    $row = $year - 1950 + 2
    $ws.Cells.Item($row, 1).Value = 266
    $ws.Cells.Item($row, 2).Value = "Gabon"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $year
    $ws.Cells.Item($row, 5).Value = "'" + $values[$year]
}
